$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 172, shifting existing rows 172:307 down to 173:308
$ws.Rows("172:172").Insert()

# Populate the newly inserted row 172 with the new record's data
$ws.Cells.Item(172, 1).Value = 4
$ws.Cells.Item(172, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(172, 3).Value = "Los Lagos"
$ws.Cells.Item(172, 4).Value = 44729
$ws.Cells.Item(172, 5).Value = 10
$ws.Cells.Item(172, 6).Value = 100114014
$ws.Cells.Item(172, 7).Value = "Betarraga"
$ws.Cells.Item(172, 8).Value = "Sin especificar"
$ws.Cells.Item(172, 9).Value = "Primera"
$ws.Cells.Item(172, 10).Value = 1000
$ws.Cells.Item(172, 11).Value = 1200
$ws.Cells.Item(172, 12).Value = 1200
$ws.Cells.Item(172, 13).Value = 1200
$ws.Cells.Item(172, 14).Value = "`$/paquete 5 unidades"
$ws.Cells.Item(172, 15).Value = "Región del Maule"
$ws.Cells.Item(172, 16).Value = 240
$ws.Cells.Item(172, 17).Value = 5
$ws.Cells.Item(172, 18).Value = "Hortaliza"
